$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record for "Repollo" (Crespo record, Primera) was added to the
# "Macroferia Regional de Talca" sheet. It belongs right above the existing
# row 290 entry (chronologically/logically), so insert a new row there,
# which pushes the former rows 290-397 down to 291-398.
$ws.Rows("290:290").Insert()

# Populate the newly inserted row 290 with the new record's data.
$ws.Range("A290").Value = 5
$ws.Range("B290").Value = "Macroferia Regional de Talca"
$ws.Range("C290").Value = "Maule"
$ws.Range("D290").Value = 44875
$ws.Range("E290").Value = 7
$ws.Range("F290").Value = 100112006
$ws.Range("G290").Value = "Repollo"
$ws.Range("H290").Value = "Crespo record"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 5000
$ws.Range("K290").Value = 1500
$ws.Range("L290").Value = 1500
$ws.Range("M290").Value = 1500
$ws.Range("N290").Value = "`$/unidad"
$ws.Range("O290").Value = "Provincia del Elquí"
$ws.Range("P290").Value = 1500
$ws.Range("Q290").Value = 1
$ws.Range("R290").Value = "Hortaliza"
